$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @("2026-01-01 10:12:43", "Admin", "Login", "login_success", "Role: admin"),
    @("2026-01-01 10:12:43", "Admin", "dashboard", "access_granted", "Opened dashboard page"),
    @("2026-01-01 10:12:46", "Admin", "quotation", "access_granted", "Opened quotation page")
)

$startRow = 208
for ($i = 0; $i -lt $newRows.Length; $i++) {
    $r = $startRow + $i
    $rowData = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $rowData[0]
    $ws.Cells.Item($r, 2).Value = $rowData[1]
    $ws.Cells.Item($r, 3).Value = $rowData[2]
    $ws.Cells.Item($r, 4).Value = $rowData[3]
    $ws.Cells.Item($r, 5).Value = $rowData[4]
}
